$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target range is formatted as Text so values remain literal strings
# (matching the workbooks original inlineStr storage) rather than being
# reinterpreted by Excel as numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("B2").Value = "14,454.00" ; $ws.Range("C2").Value = "20,462.00" ; $ws.Range("D2").Value = "17,900.56" ; $ws.Range("E2").Value = "895,028.00"
$ws.Range("B3").Value = "15,727.00" ; $ws.Range("C3").Value = "20,905.00" ; $ws.Range("D3").Value = "18,182.32" ; $ws.Range("E3").Value = "909,116.00"
$ws.Range("B4").Value = "15,148.00" ; $ws.Range("C4").Value = "20,905.00" ; $ws.Range("D4").Value = "18,611.42" ; $ws.Range("E4").Value = "930,571.00"
$ws.Range("B5").Value = "15,502.00" ; $ws.Range("C5").Value = "20,905.00" ; $ws.Range("D5").Value = "18,966.84" ; $ws.Range("E5").Value = "948,342.00"
$ws.Range("B6").Value = "16,785.00" ; $ws.Range("C6").Value = "20,905.00" ; $ws.Range("D6").Value = "19,019.22" ; $ws.Range("E6").Value = "950,961.00"
$ws.Range("B7").Value = "16,600.00" ; $ws.Range("C7").Value = "21,024.00" ; $ws.Range("D7").Value = "19,293.82" ; $ws.Range("E7").Value = "964,691.00"
$ws.Range("B8").Value = "16,321.00" ; $ws.Range("C8").Value = "21,024.00" ; $ws.Range("D8").Value = "19,703.30" ; $ws.Range("E8").Value = "985,165.00"
$ws.Range("B9").Value = "18,008.00" ; $ws.Range("C9").Value = "21,816.00" ; $ws.Range("D9").Value = "19,975.92" ; $ws.Range("E9").Value = "998,796.00"
$ws.Range("B10").Value = "18,421.00" ; $ws.Range("C10").Value = "22,077.00" ; $ws.Range("D10").Value = "20,139.84" ; $ws.Range("E10").Value = "1,006,992.00"
$ws.Range("B11").Value = "18,055.00" ; $ws.Range("C11").Value = "22,478.00" ; $ws.Range("D11").Value = "19,962.76" ; $ws.Range("E11").Value = "998,138.00"
$ws.Range("B12").Value = "17,726.00" ; $ws.Range("C12").Value = "22,478.00" ; $ws.Range("D12").Value = "20,310.52" ; $ws.Range("E12").Value = "1,015,526.00"
$ws.Range("B13").Value = "16,734.00" ; $ws.Range("C13").Value = "22,478.00" ; $ws.Range("D13").Value = "20,388.18" ; $ws.Range("E13").Value = "1,019,409.00"
$ws.Range("B14").Value = "17,135.00" ; $ws.Range("C14").Value = "22,478.00" ; $ws.Range("D14").Value = "20,315.88" ; $ws.Range("E14").Value = "1,015,794.00"
$ws.Range("B15").Value = "18,306.00" ; $ws.Range("C15").Value = "22,718.00" ; $ws.Range("D15").Value = "20,741.78" ; $ws.Range("E15").Value = "1,037,089.00"
$ws.Range("B16").Value = "18,188.00" ; $ws.Range("C16").Value = "22,946.00" ; $ws.Range("D16").Value = "20,882.82" ; $ws.Range("E16").Value = "1,044,141.00"
$ws.Range("B17").Value = "18,698.00" ; $ws.Range("C17").Value = "23,468.00" ; $ws.Range("D17").Value = "21,114.92" ; $ws.Range("E17").Value = "1,055,746.00"
$ws.Range("B18").Value = "19,293.00" ; $ws.Range("C18").Value = "23,468.00" ; $ws.Range("D18").Value = "21,059.70" ; $ws.Range("E18").Value = "1,052,985.00"
$ws.Range("B19").Value = "18,848.00" ; $ws.Range("C19").Value = "23,468.00" ; $ws.Range("D19").Value = "21,341.04" ; $ws.Range("E19").Value = "1,067,052.00"
$ws.Range("B20").Value = "19,750.00" ; $ws.Range("C20").Value = "23,468.00" ; $ws.Range("D20").Value = "21,666.50" ; $ws.Range("E20").Value = "1,083,325.00"
$ws.Range("B21").Value = "19,868.00" ; $ws.Range("C21").Value = "23,762.00" ; $ws.Range("D21").Value = "21,828.74" ; $ws.Range("E21").Value = "1,091,437.00"
$ws.Range("B22").Value = "19,288.00" ; $ws.Range("C22").Value = "23,762.00" ; $ws.Range("D22").Value = "21,958.68" ; $ws.Range("E22").Value = "1,097,934.00"
$ws.Range("B23").Value = "19,578.00" ; $ws.Range("C23").Value = "23,762.00" ; $ws.Range("D23").Value = "21,986.54" ; $ws.Range("E23").Value = "1,099,327.00"
$ws.Range("B24").Value = "19,344.00" ; $ws.Range("C24").Value = "24,229.00" ; $ws.Range("D24").Value = "22,007.46" ; $ws.Range("E24").Value = "1,100,373.00"
$ws.Range("B25").Value = "18,283.00" ; $ws.Range("C25").Value = "24,229.00" ; $ws.Range("D25").Value = "22,431.86" ; $ws.Range("E25").Value = "1,121,593.00"
$ws.Range("B26").Value = "20,086.00" ; $ws.Range("C26").Value = "24,229.00" ; $ws.Range("D26").Value = "22,607.20" ; $ws.Range("E26").Value = "1,130,360.00"
$ws.Range("B27").Value = "19,971.00" ; $ws.Range("C27").Value = "24,229.00" ; $ws.Range("D27").Value = "22,532.94" ; $ws.Range("E27").Value = "1,126,647.00"
$ws.Range("B28").Value = "19,690.00" ; $ws.Range("C28").Value = "24,229.00" ; $ws.Range("D28").Value = "22,467.70" ; $ws.Range("E28").Value = "1,123,385.00"
$ws.Range("B29").Value = "20,476.00" ; $ws.Range("C29").Value = "24,296.00" ; $ws.Range("D29").Value = "22,554.04" ; $ws.Range("E29").Value = "1,127,702.00"
$ws.Range("B30").Value = "20,268.00" ; $ws.Range("C30").Value = "24,296.00" ; $ws.Range("D30").Value = "22,490.06" ; $ws.Range("E30").Value = "1,124,503.00"
$ws.Range("B31").Value = "20,539.00" ; $ws.Range("C31").Value = "24,521.00" ; $ws.Range("D31").Value = "22,631.96" ; $ws.Range("E31").Value = "1,131,598.00"
$ws.Range("B32").Value = "20,314.00" ; $ws.Range("C32").Value = "24,691.00" ; $ws.Range("D32").Value = "22,888.00" ; $ws.Range("E32").Value = "1,144,400.00"
$ws.Range("B33").Value = "20,902.00" ; $ws.Range("C33").Value = "24,691.00" ; $ws.Range("D33").Value = "22,773.38" ; $ws.Range("E33").Value = "1,138,669.00"
$ws.Range("B34").Value = "20,445.00" ; $ws.Range("C34").Value = "24,691.00" ; $ws.Range("D34").Value = "22,862.64" ; $ws.Range("E34").Value = "1,143,132.00"
$ws.Range("B35").Value = "20,987.00" ; $ws.Range("C35").Value = "24,691.00" ; $ws.Range("D35").Value = "22,869.46" ; $ws.Range("E35").Value = "1,143,473.00"
$ws.Range("B36").Value = "21,214.00" ; $ws.Range("C36").Value = "24,691.00" ; $ws.Range("D36").Value = "23,245.76" ; $ws.Range("E36").Value = "1,162,288.00"
$ws.Range("B37").Value = "21,466.00" ; $ws.Range("C37").Value = "24,691.00" ; $ws.Range("D37").Value = "23,094.94" ; $ws.Range("E37").Value = "1,154,747.00"
$ws.Range("B38").Value = "21,427.00" ; $ws.Range("C38").Value = "24,691.00" ; $ws.Range("D38").Value = "23,125.72" ; $ws.Range("E38").Value = "1,156,286.00"
$ws.Range("B39").Value = "20,964.00" ; $ws.Range("C39").Value = "24,691.00" ; $ws.Range("D39").Value = "23,144.78" ; $ws.Range("E39").Value = "1,157,239.00"
$ws.Range("B40").Value = "22,091.00" ; $ws.Range("C40").Value = "24,691.00" ; $ws.Range("D40").Value = "23,261.16" ; $ws.Range("E40").Value = "1,163,058.00"
$ws.Range("B41").Value = "21,209.00" ; $ws.Range("C41").Value = "24,691.00" ; $ws.Range("D41").Value = "23,090.44" ; $ws.Range("E41").Value = "1,154,522.00"
$ws.Range("B42").Value = "21,673.00" ; $ws.Range("C42").Value = "24,815.00" ; $ws.Range("D42").Value = "23,248.28" ; $ws.Range("E42").Value = "1,162,414.00"
$ws.Range("B43").Value = "21,777.00" ; $ws.Range("C43").Value = "24,815.00" ; $ws.Range("D43").Value = "23,215.60" ; $ws.Range("E43").Value = "1,160,780.00"
$ws.Range("B44").Value = "20,698.00" ; $ws.Range("C44").Value = "24,815.00" ; $ws.Range("D44").Value = "23,360.08" ; $ws.Range("E44").Value = "1,168,004.00"
$ws.Range("B45").Value = "21,852.00" ; $ws.Range("C45").Value = "25,109.00" ; $ws.Range("D45").Value = "23,360.00" ; $ws.Range("E45").Value = "1,168,000.00"
$ws.Range("B46").Value = "22,094.00" ; $ws.Range("C46").Value = "25,109.00" ; $ws.Range("D46").Value = "23,460.06" ; $ws.Range("E46").Value = "1,173,003.00"
$ws.Range("B47").Value = "21,728.00" ; $ws.Range("C47").Value = "25,109.00" ; $ws.Range("D47").Value = "23,484.10" ; $ws.Range("E47").Value = "1,174,205.00"
$ws.Range("B48").Value = "21,358.00" ; $ws.Range("C48").Value = "25,109.00" ; $ws.Range("D48").Value = "23,546.54" ; $ws.Range("E48").Value = "1,177,327.00"
$ws.Range("B49").Value = "21,892.00" ; $ws.Range("C49").Value = "25,132.00" ; $ws.Range("D49").Value = "23,591.64" ; $ws.Range("E49").Value = "1,179,582.00"
$ws.Range("B50").Value = "20,444.00" ; $ws.Range("C50").Value = "25,132.00" ; $ws.Range("D50").Value = "23,585.04" ; $ws.Range("E50").Value = "1,179,252.00"
$ws.Range("B51").Value = "21,148.00" ; $ws.Range("C51").Value = "25,132.00" ; $ws.Range("D51").Value = "23,730.06" ; $ws.Range("E51").Value = "1,186,503.00"
